$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Row 2
$ws.Range("D2").Value = 19033
$ws.Range("E2").Value = 372
$ws.Range("F2").Value = 372
$ws.Range("G2").Value = -1444
$ws.Range("H2").Value = -1586
$ws.Range("I2").Value = -1574
$ws.Range("J2").Value = -11
$ws.Range("K2").Value = 20848
$ws.Range("L2").Value = 16369
$ws.Range("M2").Value = 4478
$ws.Range("N2").Value = 4381
$ws.Range("O2").Value = 97
$ws.Range("P2").Value = 2129
$ws.Range("Q2").Value = 522
$ws.Range("R2").Value = 2951
$ws.Range("S2").Value = -3049
$ws.Range("T2").Value = 153
$ws.Range("U2").Value = 370
$ws.Range("V2").Value = 9249
$ws.Range("W2").Value = 1.96
$ws.Range("X2").Value = -8.33
$ws.Range("Y2").Value = -30.51
$ws.Range("Z2").Value = -6.74
$ws.Range("AA2").Value = 365.53
$ws.Range("AB2").Value = 93.76000000000001
$ws.Range("AC2").Value = -3741
$ws.Range("AD2").Value = -1.74
$ws.Range("AE2").Value = 10287
$ws.Range("AF2").Value = 0.63
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = -6.03
$ws.Range("AJ2").Value = 32414942

# Row 3
$ws.Range("D3").Value = 18553
$ws.Range("E3").Value = 310
$ws.Range("F3").Value = 310
$ws.Range("G3").Value = -998
$ws.Range("H3").Value = -1144
$ws.Range("I3").Value = -1142
$ws.Range("J3").Value = -2
$ws.Range("K3").Value = 22854
$ws.Range("L3").Value = 19622
$ws.Range("M3").Value = 3233
$ws.Range("N3").Value = 3132
$ws.Range("O3").Value = 100
$ws.Range("P3").Value = 2183
$ws.Range("Q3").Value = 944
$ws.Range("R3").Value = -2302
$ws.Range("S3").Value = 1261
$ws.Range("T3").Value = 137
$ws.Range("U3").Value = 807
$ws.Range("V3").Value = 9804
$ws.Range("W3").Value = 1.67
$ws.Range("X3").Value = -6.17
$ws.Range("Y3").Value = -30.41
$ws.Range("Z3").Value = -5.24
$ws.Range("AA3").Value = 606.99
$ws.Range("AB3").Value = 38.51
$ws.Range("AC3").Value = -2635
$ws.Range("AD3").Value = -1.62
$ws.Range("AE3").Value = 7173
$ws.Range("AF3").Value = 0.6
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 33494942

# Row 4
$ws.Range("D4").Value = 18318
$ws.Range("E4").Value = 955
$ws.Range("F4").Value = 955
$ws.Range("G4").Value = 257
$ws.Range("H4").Value = 102
$ws.Range("I4").Value = 105
$ws.Range("J4").Value = -3
$ws.Range("K4").Value = 20739
$ws.Range("L4").Value = 17221
$ws.Range("M4").Value = 3518
$ws.Range("N4").Value = 3420
$ws.Range("O4").Value = 97
$ws.Range("P4").Value = 2333
$ws.Range("Q4").Value = 2080
$ws.Range("R4").Value = 572
$ws.Range("S4").Value = -3051
$ws.Range("T4").Value = 452
$ws.Range("U4").Value = 1628
$ws.Range("V4").Value = 6338
$ws.Range("W4").Value = 5.21
$ws.Range("X4").Value = 0.5600000000000001
$ws.Range("Y4").Value = 3.2
$ws.Range("Z4").Value = 0.47
$ws.Range("AA4").Value = 489.53
$ws.Range("AB4").Value = 42.1
$ws.Range("AC4").Value = 231
$ws.Range("AD4").Value = 18.44
$ws.Range("AE4").Value = 7329
$ws.Range("AF4").Value = 0.58
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 23.84
$ws.Range("AJ4").Value = 36494942

# Row 5
$ws.Range("D5").Value = 16415
$ws.Range("E5").Value = 1495
$ws.Range("F5").Value = 1495
$ws.Range("G5").Value = 753
$ws.Range("H5").Value = 466
$ws.Range("I5").Value = 464
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 22141
$ws.Range("L5").Value = 18215
$ws.Range("M5").Value = 3926
$ws.Range("N5").Value = 3827
$ws.Range("O5").Value = 100
$ws.Range("P5").Value = 2333
$ws.Range("Q5").Value = 770
$ws.Range("R5").Value = -2706
$ws.Range("S5").Value = 1536
$ws.Range("T5").Value = 1865
$ws.Range("U5").Value = -1095
$ws.Range("V5").Value = 7969
$ws.Range("W5").Value = 9.109999999999999
$ws.Range("X5").Value = 2.84
$ws.Range("Y5").Value = 12.8
$ws.Range("Z5").Value = 2.17
$ws.Range("AA5").Value = 463.91
$ws.Range("AB5").Value = 62.93
$ws.Range("AC5").Value = 994
$ws.Range("AD5").Value = 4.1
$ws.Range("AE5").Value = 8382
$ws.Range("AF5").Value = 0.49
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 2.45
$ws.Range("AI5").Value = 84.11
$ws.Range("AJ5").Value = 36494942

# Row 6
$ws.Range("D6").Value = 13210
$ws.Range("E6").Value = 602
$ws.Range("F6").Value = 602
$ws.Range("G6").Value = 92
$ws.Range("H6").Value = -121
$ws.Range("I6").Value = -124
$ws.Range("K6").Value = 17938
$ws.Range("L6").Value = 15281
$ws.Range("M6").Value = 2657
$ws.Range("N6").Value = 2554
$ws.Range("P6").Value = 2452
$ws.Range("Q6").Value = 1337
$ws.Range("R6").Value = 3426
$ws.Range("S6").Value = -4218
$ws.Range("T6").Value = 2429
$ws.Range("U6").Value = -1091
$ws.Range("V6").Value = 3716
$ws.Range("W6").Value = 4.55
$ws.Range("X6").Value = -0.91
$ws.Range("Y6").Value = -3.89
$ws.Range("Z6").Value = -0.6
$ws.Range("AA6").Value = 575.2
$ws.Range("AB6").Value = 4.71
$ws.Range("AC6").Value = -259
$ws.Range("AD6").Value = -16.63
$ws.Range("AE6").Value = 5316
$ws.Range("AF6").Value = 0.8100000000000001
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 38874942

# Row 7
$ws.Range("D7").Value = 13236
$ws.Range("E7").Value = 544
$ws.Range("G7").Value = 144
$ws.Range("H7").Value = 38
$ws.Range("I7").Value = 30
$ws.Range("K7").Value = 17660
$ws.Range("L7").Value = 14978
$ws.Range("M7").Value = 2682
$ws.Range("N7").Value = 2545
$ws.Range("P7").Value = 2452
$ws.Range("Q7").Value = -263
$ws.Range("R7").Value = 372
$ws.Range("S7").Value = -159
$ws.Range("T7").Value = 644
$ws.Range("U7").Value = 864
$ws.Range("W7").Value = 4.11
$ws.Range("X7").Value = 0.29
$ws.Range("Y7").Value = 1.16
$ws.Range("Z7").Value = 0.22
$ws.Range("AA7").Value = 558.46
$ws.Range("AC7").Value = 60
$ws.Range("AD7").Value = 44.73
$ws.Range("AE7").Value = 5298
$ws.Range("AF7").Value = 0.51
$ws.Range("AG7").Value = 50
$ws.Range("AH7").Value = 1.86
$ws.Range("AI7").Value = 65.89

# Row 8
$ws.Range("D8").Value = 13924
$ws.Range("E8").Value = 671
$ws.Range("G8").Value = 364
$ws.Range("H8").Value = 236
$ws.Range("I8").Value = 218
$ws.Range("K8").Value = 18639
$ws.Range("L8").Value = 15788
$ws.Range("M8").Value = 2852
$ws.Range("N8").Value = 2581
$ws.Range("P8").Value = 2452
$ws.Range("Q8").Value = 799
$ws.Range("R8").Value = -830
$ws.Range("S8").Value = -122
$ws.Range("T8").Value = 739
$ws.Range("U8").Value = 420
$ws.Range("W8").Value = 4.82
$ws.Range("X8").Value = 1.7
$ws.Range("Y8").Value = 8.49
$ws.Range("Z8").Value = 1.3
$ws.Range("AA8").Value = 553.66
$ws.Range("AC8").Value = 443
$ws.Range("AD8").Value = 6.07
$ws.Range("AE8").Value = 5373
$ws.Range("AF8").Value = 0.5
$ws.Range("AG8").Value = 50
$ws.Range("AH8").Value = 1.86
$ws.Range("AI8").Value = 8.94

# Row 9
$ws.Range("D9").Value = 14752
$ws.Range("E9").Value = 806
$ws.Range("G9").Value = 522
$ws.Range("H9").Value = 356
$ws.Range("I9").Value = 328
$ws.Range("K9").Value = 19496
$ws.Range("L9").Value = 16336
$ws.Range("M9").Value = 3159
$ws.Range("N9").Value = 2737
$ws.Range("P9").Value = 2452
$ws.Range("Q9").Value = 1010
$ws.Range("R9").Value = -823
$ws.Range("S9").Value = -102
$ws.Range("T9").Value = 712
$ws.Range("U9").Value = 698
$ws.Range("W9").Value = 5.47
$ws.Range("X9").Value = 2.41
$ws.Range("Y9").Value = 12.35
$ws.Range("Z9").Value = 1.86
$ws.Range("AA9").Value = 517.14
$ws.Range("AC9").Value = 670
$ws.Range("AD9").Value = 4.02
$ws.Range("AE9").Value = 5698
$ws.Range("AF9").Value = 0.47
$ws.Range("AG9").Value = 50
$ws.Range("AH9").Value = 1.86
$ws.Range("AI9").Value = 5.92
